$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 3-18 (row 18 being "Patient 15").
# Add a new row 19 ("Patient 16") describing a patient sharing an NHS number
# with "Patient3", mirroring the formatting of the row above it.

$ws.Range("A18:O18").Copy()
$ws.Range("A19:O19").PasteSpecial(-4122)

$ws.Range("A19").Value = "Patient 16"
$ws.Range("B19").Value = "?"
$ws.Range("C19").Value = "?"
$ws.Range("D19").Value = "?"
$ws.Range("E19").Value = "?"
$ws.Range("F19").Value = "?"
$ws.Range("G19").Value = "?"
$ws.Range("H19").Value = "?"
$ws.Range("I19").Value = "?"
$ws.Range("J19").Value = "?"
$ws.Range("K19").Value = "?"
$ws.Range("L19").Value = "?"
$ws.Range("M19").Value = "?"
$ws.Range("N19").Value = "?"
$ws.Range("O19").Value = "Patient with same NHS number as “Patient3”"

# Grey out the "don't care" marks for the new scenario row so it reads
# differently to the regular "x"/"?" rows above it.
$ws.Range("B19:N19").Font.Color = 11711154

$ws.Range("D22").Select()
